# "Generate Report for Handoff"
#
# Replaces the two pending localization entries (734d1017-...md /
# c20ec0b6-...md) with a fresh handoff pair (0986ad25-...md /
# ffffdaf03d5b-...md), flips the status from
# "Handed back: in sync with en-US" to "Ready for handoff", refreshes the
# handoff file names + timestamps, and drops the now-obsolete
# "Latest Target File" / "Latest Handback File" columns (the handback
# hasn't happened yet) on the per-locale sheets.

$wb = $excel.ActiveWorkbook

# ---- old / new identifiers -------------------------------------------------
$newMd1 = "0986ad25-7ee2-4c21-9928-79827311773c.md"
$newMd2 = "ffffdaf03d5b-b8c2-41fe-9a94-58525e86c329.md"
$newStatus = "Ready for handoff"
$pendingDt = "0001-01-01 00:00:00"

$newXlfZh = "0986ad25-7ee2-4c21-9928-79827311773c.32803ce57d513cf9f6d2829fe3358f10787d7b53.zh-cn.xlf"
$newXlfDe = "0986ad25-7ee2-4c21-9928-79827311773c.32803ce57d513cf9f6d2829fe3358f10787d7b53.de-de.xlf"

$newHandoffDtZh = "2016-03-08 12:47:50"
$newHandoffDtDe = "2016-03-08 12:47:54"

$mdBaseUrl = "https://github.com/OpenLocalizationTest/oltest/blob/fbb0154626b8c336e50968d3e4d1f3278f3523f6/e2e/"

# ---- Overview sheet ---------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# NOTE: in this COM host, Range(...).Hyperlinks.Delete() clears *every*
# hyperlink on the sheet (not just the ones intersecting that range), and
# the single-Hyperlink .Delete() method is a no-op. So: wipe the sheet's
# hyperlinks exactly once, then re-Add every hyperlink the sheet should end
# up with (changed ones with new data, untouched ones with their original
# data).
$wsOverview.Range("A1").Hyperlinks.Delete()

$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), ($mdBaseUrl + $newMd1), "", "", $newMd1)
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), ($mdBaseUrl + $newMd2), "", "", $newMd2)
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), ($mdBaseUrl + ".localization-config"), "", "", ".localization-config")

$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

# ---- per-locale sheets -------------------------------------------------------
function Update-LocaleSheet($ws, $newMd1, $newMd2, $newStatus, $newXlf, $newHandoffDt, $pendingDt, $mdBaseUrl, $handoffBaseUrl) {
    # Wipe this sheet's hyperlinks exactly once, up front.
    $ws.Range("A1").Hyperlinks.Delete()

    # Row 2 -----------------------------------------------------------------
    $ws.Hyperlinks.Add($ws.Range("A2"), ($mdBaseUrl + $newMd1), "", "", $newMd1)
    $ws.Range("B2").Value = $newStatus

    $ws.Hyperlinks.Add($ws.Range("C2"), ($handoffBaseUrl + $newXlf), "", "", $newXlf)
    $ws.Range("D2").Value = $newHandoffDt

    $ws.Range("E2").Clear()
    $ws.Range("F2").Clear()

    $ws.Range("G2").Value = $pendingDt
    $ws.Range("H2").Value = "Include"

    # Row 3 -----------------------------------------------------------------
    $ws.Hyperlinks.Add($ws.Range("A3"), ($mdBaseUrl + $newMd2), "", "", $newMd2)
    $ws.Range("B3").Value = $newStatus

    $ws.Hyperlinks.Add($ws.Range("C3"), ($handoffBaseUrl + $newXlf), "", "", $newXlf)
    $ws.Range("D3").Value = $newHandoffDt

    $ws.Range("E3").Clear()
    $ws.Range("F3").Clear()

    $ws.Range("G3").Value = $pendingDt
    $ws.Range("H3").Value = "Include"

    # Row 4 (".localization-config" / "Not to be localized") is otherwise
    # untouched, but its hyperlink needs to be re-added since the sheet's
    # hyperlink collection was wiped above.
    $ws.Hyperlinks.Add($ws.Range("A4"), ($mdBaseUrl + ".localization-config"), "", "", ".localization-config")
}

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$zhHandoffBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6128a206767ab4b5b6f3c3792f535b483b61834b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/"
Update-LocaleSheet $wsZhCn $newMd1 $newMd2 $newStatus $newXlfZh $newHandoffDtZh $pendingDt $mdBaseUrl $zhHandoffBase

$wsDeDe = $wb.Worksheets.Item("de-de")
$deHandoffBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0b71a660de2c0a8ff2662dd5f73e6a8c05ba5ced/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/"
Update-LocaleSheet $wsDeDe $newMd1 $newMd2 $newStatus $newXlfDe $newHandoffDtDe $pendingDt $mdBaseUrl $deHandoffBase
